$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.195854663848877
$ws.Range("B1").Value = 2.595503568649292
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.188163757324219
$ws.Range("E1").Value = 1.178524613380432
